# moving frames in folder corresponding to class
# Adds a new "classes" worksheet (after "labels") that maps the numeric
# class id used in the "labels" sheet ("class" column) to its textual
# class name.

$wb = $excel.ActiveWorkbook

$delayFrames = $wb.Worksheets.Item("delay_frames")
$labels = $wb.Worksheets.Item("labels")

# --- Insert the new "classes" sheet right after "labels" -------------------
$classes = $wb.Worksheets.Add($null, $labels)
$classes.Name = "classes"

# --- Populate the class lookup table ----------------------------------------
# Header
$classes.Range("B1").Value = "class"

# Write the B-column (text) values first, in the same order the class names
# were first authored (alphabetical-ish "definition" order), so the shared
# string table is built in that exact sequence - then fill in row 2 last.
$classes.Range("B3").Value = "falling"
$classes.Range("B4").Value = "lying_on_the_ground"
$classes.Range("B5").Value = "crounching"
$classes.Range("B6").Value = "moving_down"
$classes.Range("B7").Value = "moving_up"
$classes.Range("B8").Value = "sitting"
$classes.Range("B9").Value = "lying_on_a_sofa"
$classes.Range("B10").Value = "moving_horizontaly"
$classes.Range("B2").Value = "walking_or_standing_up"

# Id column
$classes.Range("A2").Value = 1
$classes.Range("A3").Value = 2
$classes.Range("A4").Value = 3
$classes.Range("A5").Value = 4
$classes.Range("A6").Value = 5
$classes.Range("A7").Value = 6
$classes.Range("A8").Value = 7
$classes.Range("A9").Value = 8
$classes.Range("A10").Value = 9

# --- Column widths -----------------------------------------------------------
$classes.Columns.Item(1).ColumnWidth = 24.166666666666668   # -> stored width 25
$classes.Columns.Item(2).ColumnWidth = 17.666666666666668   # -> stored width 18.5

# --- View / selection state --------------------------------------------------
# delay_frames: selection moves from J24 to J34 (no tab activation)
$delayFrames.Range("J34").Select()

# labels: selection moves from F152 to G164, it is no longer the
# scrolled/top-left-pinned, tab-selected sheet
$labels.Range("G164").Select()

# classes becomes the active / tab-selected sheet, zoomed to 243%, with B2
# selected
$classes.Activate()
$classes.Range("B2").Select()
$excel.ActiveWindow.Zoom = 243
